$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "64.170.79"
$ws.Range("E2").Value = "  -0.85%  "
# Row 3
$ws.Range("D3").Value = "3.402.13"
$ws.Range("E3").Value = "  -1.07%  "
# Row 4
$ws.Range("E4").Value = "  +0.00%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "569.30"
$ws.Range("E5").Value = "  -1.10%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "155.39"
$ws.Range("E6").Value = "  -2.98%  "
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.624"
$ws.Range("E7").Value = "  +7.37%  "
# Row 8
$ws.Range("E8").Value = "  +0.01%  "
# Row 9
$ws.Range("D9").Value = "3.405.80"
$ws.Range("E9").Value = "  -1.07%  "
# Row 10
$ws.Range("E10").Value = "  -3.03%  "
# Row 11
$ws.Range("E11").Value = "  -2.62%  "
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.438"
$ws.Range("E12").Value = "  -0.69%  "
# Row 13
$ws.Range("D13").Value = "3.986.82"
$ws.Range("E13").Value = "  -1.21%  "
# Row 14
$ws.Range("E14").Value = "  -0.16%  "
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000186"
$ws.Range("E15").Value = "  -3.72%  "
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "27.60"
$ws.Range("E16").Value = "  -2.34%  "
# Row 17
$ws.Range("D17").Value = "64.180.93"
$ws.Range("E17").Value = "  -0.85%  "
# Row 18
$ws.Range("D18").Value = "3.403.82"
$ws.Range("E18").Value = "  -0.75%  "
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.29"
$ws.Range("E19").Value = "  -1.20%  "
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.81"
$ws.Range("E20").Value = "  -3.30%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "375.56"
$ws.Range("E21").Value = "  -2.82%  "
# Row 22
$ws.Range("E22").Value = "  -2.31%  "
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.544"
$ws.Range("E23").Value = "  -0.15%  "
# Row 24
$ws.Range("E24").Value = "  -0.21%  "
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "71.55"
$ws.Range("E25").Value = "  -2.31%  "
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000118"
$ws.Range("E26").Value = "  -4.51%  "
# Row 27
$ws.Range("E27").Value = "  +6.22%  "
# Row 28
$ws.Range("E28").Value = "  -1.69%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.01"
$ws.Range("E29").Value = "  +0.55%  "
# Row 30
$ws.Range("E30").Value = "  +2.45%  "
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.16"
$ws.Range("E31").Value = "  -0.23%  "
# Row 32
$ws.Range("E32").Value = "  -2.16%  "
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.01"
$ws.Range("E33").Value = "  -2.71%  "
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.12"
$ws.Range("E34").Value = "  +0.44%  "
# Row 35
$ws.Range("E35").Value = "  +5.80%  "
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "159.53"
$ws.Range("E36").Value = "  -2.25%  "
# Row 37
$ws.Range("E37").Value = "  -0.17%  "
# Row 38
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0758"
$ws.Range("E38").Value = "  -0.94%  "
# Row 39
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.85"
$ws.Range("E39").Value = "  +4.45%  "
# Row 40
$ws.Range("D40").Value = "2.872.79"
$ws.Range("E40").Value = "  -4.95%  "
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "26.27"
$ws.Range("E41").Value = "  -3.52%  "
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.59"
$ws.Range("E42").Value = "  +1.00%  "
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "42.81"
$ws.Range("E43").Value = "  -0.03%  "
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "26.25"
$ws.Range("E44").Value = "  +6.34%  "
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0314"
$ws.Range("E45").Value = "  -0.79%  "
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.767"
$ws.Range("E46").Value = "  -0.70%  "
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "319.70"
$ws.Range("E47").Value = "  +6.03%  "
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.07"
$ws.Range("E48").Value = "  -1.35%  "
# Row 49
$ws.Range("E49").Value = "  +2.49%  "
# Row 50
$ws.Range("E50").Value = "  +0.87%  "
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.54"
$ws.Range("E51").Value = "  -1.32%  "
